$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the standalone "Meta description: ..." paragraph that currently
#    sits right after the H1 title.
# ---------------------------------------------------------------------------
$metaDeleted = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Meta description:*") {
        $p.Range.Delete()
        $metaDeleted = $true
        break
    }
}

# ---------------------------------------------------------------------------
# 2) Insert a new bold paragraph "Play Dragon Tribe Slot Game for Free -
#    Exciting Dragon Spins Feature" right before the final paragraph (the
#    italic "image prompt" paragraph).
# ---------------------------------------------------------------------------
$oldCount = $d.Paragraphs.Count
$targetIdx = $oldCount
for ($i = 1; $i -le $oldCount; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "Create a feature image for Dragon Tribe*") {
        $targetIdx = $i
        break
    }
}
$lastPara = $d.Paragraphs.Item($targetIdx)
$insertPos = $lastPara.Range.Start
$insertRange = $d.Range($insertPos, $insertPos)

$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Dragon Tribe Slot Game for Free - Exciting Dragon Spins Feature</w:t></w:r></w:p><w:p/></w:body></w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@

$insertRange.InsertXML($xml)

# InsertXML above creates two new paragraphs (the bold one + a throwaway
# empty one) to force a real paragraph break; collapse the stray empty
# paragraph back out so the bold paragraph sits directly before the final
# (image prompt) paragraph, with no blank paragraph between them.
$strayIdx = $targetIdx + 1
$strayPara = $d.Paragraphs.Item($strayIdx)
if ($strayPara.Range.Text.Trim() -eq "") {
    $strayPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 3) Replace the final paragraph's text (the old "Create a feature image..."
#    AI-art prompt) with the new meta-description copy, keeping its existing
#    (italic) run formatting.
# ---------------------------------------------------------------------------
$oldText = "Create a feature image for Dragon Tribe, a cartoon-style image featuring a happy Maya warrior with glasses. The image should incorporate the dragon theme with flames and dragons flying in the background. Use bold and bright colors to make the image pop and catch the eye of potential players. Make sure to showcase the xNudge Wilds and Dragon Spins features in the image to give players a taste of the action-packed gameplay. The Maya warrior should be standing in front of the reel set, with the game's logo at the top and the words ""Dragon Tribe"" written in a fun and playful font. Overall, the image should convey the excitement and adventure of this slot game."
$newText = "Experience the dragon-filled world of Dragon Tribe slots. Play for free and win up to 27,000x with the exciting Dragon Spins feature and xNudge Wilds."

$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
